$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (number formats, borders, fonts) from row 7 down to row 8
# so the new row matches the existing table style exactly.
$ws.Range("B7:N7").Copy()
$ws.Range("B8:N8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the values for the new row (row 8)
# (columns D, E, F, G, I are intentionally left blank, same as row 7's I:L before data)
$ws.Range("B8").Value = "D2_backup"
$ws.Range("C8").Value = "D2_backup_dict"
$ws.Range("H8").Value = "mo_D2_backup"
$ws.Range("J8").Value = 1134
$ws.Range("K8").Value = 1134
$ws.Range("L8").Value = 1134
$ws.Range("M8").Value = 43339
$ws.Range("N8").Formula = "=M8-M7"

# Setting a date-difference formula in N8 makes Excel auto-apply a date
# number format to the cell; restore the original (General/no-border-change)
# formatting by re-copying it from N7.
$ws.Range("N7").Copy()
$ws.Range("N8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the sheet view: active selection moves to the new row's K cell
$ws.Range("K8").Select()
